$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows that are fully removed in the target layout
$ws.Range("A33:I33").ClearContents()
$ws.Range("A37:I38").ClearContents()

# Cell value updates
$data = @(
    ,@("C8", "num", 29)
    ,@("C9", "num", 79)
    ,@("G9", "text", '37288.00')
    ,@("C10", "num", 13)
    ,@("G10", "text", '8606.00')
    ,@("A11", "text", 'P. point')
    ,@("C11", "num", 54)
    ,@("D11", "text", '6')
    ,@("E11", "text", 'On board')
    ,@("F11", "num", 136)
    ,@("G11", "text", '7344.00')
    ,@("A12", "text", 'Each')
    ,@("C12", "num", 56)
    ,@("D12", "text", '3.0')
    ,@("E12", "text", 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .')
    ,@("F12", "num", 23)
    ,@("G12", "text", '1288.00')
    ,@("C13", "num", 58)
    ,@("D13", "text", '4.0')
    ,@("E13", "text", 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .')
    ,@("F13", "num", 50)
    ,@("G13", "text", '2900.00')
    ,@("C14", "num", 10)
    ,@("D14", "text", '6.0')
    ,@("E14", "text", 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .')
    ,@("F14", "num", 78)
    ,@("G14", "text", '780.00')
    ,@("C15", "num", 2)
    ,@("D15", "text", '7.0')
    ,@("E15", "text", 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .')
    ,@("F15", "num", 30)
    ,@("G15", "text", '60.00')
    ,@("C16", "num", 75)
    ,@("D16", "text", '8.0')
    ,@("E16", "text", 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .')
    ,@("F16", "num", 30)
    ,@("G16", "text", '2250.00')
    ,@("C17", "num", 55)
    ,@("D17", "text", '9.0')
    ,@("E17", "text", 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .')
    ,@("F17", "num", 219)
    ,@("G17", "text", '12045.00')
    ,@("A18", "text", 'Each')
    ,@("C18", "num", 52)
    ,@("D18", "text", '10.0')
    ,@("E18", "text", 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR')
    ,@("F18", "num", 303)
    ,@("G18", "text", '15756.00')
    ,@("C19", "num", 41)
    ,@("G19", "text", '1640.00')
    ,@("A20", "text", '')
    ,@("C20", "num", 49)
    ,@("D20", "text", '12.0')
    ,@("E20", "text", 'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR')
    ,@("F20", "num", 0)
    ,@("G20", "text", '0.00')
    ,@("A21", "text", 'Mtr.')
    ,@("C21", "num", 5)
    ,@("D21", "text", '20')
    ,@("E21", "text", '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.')
    ,@("F21", "num", 122)
    ,@("G21", "text", '610.00')
    ,@("C22", "num", 78)
    ,@("C23", "num", 26)
    ,@("G23", "text", '520.00')
    ,@("C24", "num", 80)
    ,@("C25", "num", 58)
    ,@("G25", "text", '109620.00')
    ,@("C26", "num", 88)
    ,@("D26", "text", '17.0')
    ,@("E26", "text", 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR')
    ,@("A27", "text", '')
    ,@("C27", "num", 6)
    ,@("D27", "text", '29')
    ,@("E27", "text", 'Single pole MCB   (With B/C curve tripping Characteristics)')
    ,@("F27", "num", 0)
    ,@("G27", "text", '0.00')
    ,@("D28", "text", '18.0')
    ,@("E28", "text", 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR')
    ,@("A29", "text", '')
    ,@("C29", "num", 80)
    ,@("D29", "text", '34')
    ,@("E29", "text", 'Metal door (single phase) IK-09 and IP-43 with Metal end box')
    ,@("F29", "num", 0)
    ,@("G29", "text", '0.00')
    ,@("A30", "text", 'Each')
    ,@("C30", "num", 6)
    ,@("D30", "text", '35')
    ,@("E30", "text", '8 Way (8+2)')
    ,@("F30", "num", 2184)
    ,@("G30", "text", '13104.00')
    ,@("C31", "num", 62)
    ,@("D31", "text", '36')
    ,@("E31", "text", 'Total')
    ,@("A32", "text", '')
    ,@("C32", "num", 98)
    ,@("D32", "text", '38')
    ,@("E32", "text", 'Grand Total')
    ,@("F32", "num", 0)
    ,@("G32", "text", '0.00')
    ,@("B34", "text", '')
    ,@("C34", "text", '')
    ,@("D34", "text", '')
    ,@("E34", "text", 'Grand Total Rs.')
    ,@("F34", "text", '')
    ,@("G34", "text", '213811.00')
    ,@("H34", "text", '213811.00')
    ,@("I34", "text", '')
    ,@("A35", "text", '')
    ,@("B35", "text", '')
    ,@("C35", "text", '')
    ,@("D35", "text", '')
    ,@("E35", "text", 'Tender Premium @ 0%')
    ,@("F35", "text", '')
    ,@("G35", "text", '0.00')
    ,@("H35", "text", '0.00')
    ,@("I35", "text", '')
    ,@("A36", "text", '')
    ,@("B36", "text", '')
    ,@("C36", "text", '')
    ,@("D36", "text", '')
    ,@("E36", "text", 'NET PAYABLE AMOUNT Rs.')
    ,@("F36", "text", '')
    ,@("G36", "text", '213811.00')
    ,@("H36", "text", '213811.00')
    ,@("I36", "text", '')
)

foreach ($row in $data) {
    $cellRef = $row[0]
    $kind = $row[1]
    $val = $row[2]
    if ($kind -eq "text") {
        $ws.Range($cellRef).Value = "'" + $val
    } elseif ($kind -eq "clear") {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $val
    }
}
